# Update the "取得日時" (retrieved timestamp) column for every data row on
# the "ランサーズ" sheet to reflect the latest scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-08 06:34:41"

# Data rows run from row 2 through row 11 (row 1 is the header row).
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
